# "Changes for the next release"
# Updates a handful of input cells on Sheet1 (columns E/F for rows 13-15 of
# the 启蒙/工业/现代 output-balance table). The H (=SUM(B:G)) and J:O
# (period-over-period delta) columns are driven by shared formulas, so they
# recalculate automatically once the underlying inputs change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 工业 row (13): F 1.6 -> 1.55
$ws.Range("F13").Value = 1.55

# 现代（电气以后） row (14): E 1.5 -> 1.6, F 1.6 -> 1.55
$ws.Range("E14").Value = 1.6
$ws.Range("F14").Value = 1.55

# row (15): E 1.5 -> 1.6, F 1.6 -> 1.55
$ws.Range("E15").Value = 1.6
$ws.Range("F15").Value = 1.55

# Leave the selection where the author ended up after the edit.
[void]$ws.Range("N12:N13").Select()
